# Scrum board sprint 3 update
# - add a "Backlog" / "Digital in" note pair in column I next to sprint 2 rows
# - move several DOING (column E) checkmarks to their proper TODO/DONE column
# - fill in sprint 3 (rows 15-19) task description / owner / TODO marks

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New notes in column I (next to row 7 / row 8) ---
$ws.Range("I7").Value = "Backlog"
$ws.Range("I7").Font.Bold = $true

$ws.Range("I8").Value = "Digital in"

# --- Sprint 2 status corrections (checkmark moved between TODO/DOING/DONE columns) ---
# Row 9 (2.1 Kernel module): was marked DOING (E9=X) -> now DONE (F9=X)
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = "X"

# Row 11 (2.3 Digital in): was marked DOING (E11=X) -> now TODO (D11=X)
$ws.Range("D11").Value = "X"
$ws.Range("E11").Value = ""

# Row 12 (2.4 PWM): was marked DOING (E12=X) -> now DONE (F12=X)
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = "X"

# Row 13 (2.5 I2C): was marked DOING (E13=X) -> now DONE (F13=X)
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = "X"

# --- Sprint 3 backlog items (rows 15-19) ---
# Row 15 (3.1)
$ws.Range("B15").Value = "cross compile module for pi"
$ws.Range("C15").Value = "Everyone"
$ws.Range("D15").Value = "X"

# Row 16 (3.2)
$ws.Range("B16").Value = "blink led module for pi"
$ws.Range("C16").Value = "Bart"
$ws.Range("D16").Value = "X"

# Row 17 (3.3)
$ws.Range("B17").Value = "porting PWM naar pi2"
$ws.Range("C17").Value = "Dennis/Robbert"
$ws.Range("D17").Value = "X"

# Row 18 (3.4)
$ws.Range("B18").Value = "algoritme besturing"
$ws.Range("C18").Value = "Jeroen"
$ws.Range("D18").Value = "X"

# Row 19 (3.5)
$ws.Range("B19").Value = "software architectuur"
$ws.Range("C19").Value = "Dennis"
$ws.Range("D19").Value = "X"
